$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value, $donorRef) {
    $ws.Range($cellRef).Value = "'" + $value
    $ws.Range($donorRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)  # xlPasteFormats
}

Set-TextValue "E2" "15" "B9"
Set-TextValue "E3" "21" "B9"
Set-TextValue "B12" "44" "B9"
Set-TextValue "C12" "15" "B9"
Set-TextValue "B13" "17" "B9"
Set-TextValue "C13" "44" "B9"
Set-TextValue "B14" "20" "B9"
Set-TextValue "C14" "22" "B9"
Set-TextValue "C16" "20" "B9"
Set-TextValue "D22" "18" "B9"
Set-TextValue "D23" "22" "B9"
Set-TextValue "D24" "44" "B9"
Set-TextValue "D25" "20" "B9"
